# Convert M2Doc field-code paragraphs ("begin"/instrText/"separate"/"end")
# into plain-text runs wrapped in curly braces, e.g.
#   {begin}{instrText "m:for i | Sequence{1, 2, 3}"}{sep}...{end}
# becomes four runs of plain text:
#   "{m:for i | Sequence{1, 2, 3}" / "" / "" / "}"
# This mirrors TokenIteratorFieldRewriterSplit turning field codes back
# into literal `{...}` token text.

function Escape-XmlText([string]$s) {
    if ($null -eq $s) { return "" }
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

$d = $word.ActiveDocument

while ($d.Fields.Count -gt 0) {
    $f = $d.Fields.Item(1)
    $code = $f.Code.Text
    $codeStart = $f.Code.Start

    # Locate the paragraph that contains this field's instruction text.
    $target = $null
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($codeStart -ge $p.Range.Start -and $codeStart -lt $p.Range.End) {
            $target = $p
            break
        }
    }

    $openText = Escape-XmlText ("{" + $code)

    $newBody = "<w:p>" +
        "<w:r><w:rPr/><w:t>" + $openText + "</w:t></w:r>" +
        "<w:r><w:rPr/><w:t/></w:r>" +
        "<w:r><w:rPr/><w:t/></w:r>" +
        "<w:r><w:rPr/><w:t>}</w:t></w:r>" +
        "</w:p>"

    $xml = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $newBody +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # Insert the replacement runs at the (collapsed) start of the
    # paragraph - this leaves the paragraph mark / pPr completely
    # untouched (InsertXML on a non-collapsed range that spans the
    # paragraph mark rebuilds pPr and drops empty <w:rPr/> containers).
    $insPoint = $d.Range($target.Range.Start, $target.Range.Start)
    $insPoint.InsertXML($xml)

    # Now remove the original field (begin/instrText/separate/.../end
    # runs); it got pushed later in the paragraph by the insert above.
    $f.Delete()
}

Write-Output "done"
